# 5.17 Add Character Assets
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change F2 from "Waterfall" to "Happy" (adds a new shared string "Happy"
# and removes the now-unused "Waterfall" string from the shared strings table).
# Dependent formulas in N3:N6 (=IF(F2<>"",F2,N2) chain) will recalc to "Happy".
$ws.Range("F2").Value = "Happy"

# Update the active selection on the sheet view from B16 to F4.
$ws.Range("F4").Select()
